$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "leadlag")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 7

    $ws.Cells.Item($row, 1).Value = 6

    # Force the date column to stay a plain text string instead of being
    # auto-converted into a date serial value by Excel's input parsing.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-16"

    $ws.Cells.Item($row, 3).Value = "21:51:32"
    $ws.Cells.Item($row, 4).Value = "leadlag"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 68263.935
    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100
    $ws.Cells.Item($row, 12).Value = 0.6267
    $ws.Cells.Item($row, 13).Value = "Binance leading with -0.063% move"
    $ws.Cells.Item($row, 15).Value = 0
}
